$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row, columns B..G
$data = @{
    2  = @(1.445647641019636, 0.3048912486333797, 18.71679738969934, 13.86384647080068, 0, 34.33118275015303)
    3  = @(0.6545652718822623, 0.3048912486333797, 0.7210945179870265, 0.5333859586016987, 1, 2.213936997104367)
    4  = @(1.445647641019636, 1.626987699542094, 189.6080260415259, 0.5333859586016987, 0, 193.2140473406893)
    5  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 0, 8.656069925401464)
    6  = @(0.04172184405617529, 0.3048912486333797, 18.71679738969934, 0.5333859586016987, 1, 19.59679644099059)
    7  = @(0.6545652718822623, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 6.038307959104277)
    8  = @(0.0006075818656279264, 0.002658071450198252, 0.7210945179870265, 0.5333859586016987, 1, 1.257746129904551)
    9  = @(3.272327238179451, 1.626987699542094, 3.223369029078222, 0.5333859586016987, 1, 8.656069925401464)
    10 = @(3.272327238179451, 1.626987699542094, 18.71679738969934, 0.5333859586016987, 1, 24.14949828602258)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
